$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New numeric column D: "Predictor" simulated probability, formatted as 0.00% ---
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0.5
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("D3:D8").NumberFormat = "0.00%"

# --- New comment column E ---
$ws.Range("E3").Value = "Main loop is always taken(unconditional jump) and we have 50% of all predictions as mispredicted. Tested branch (compare 0 and 0) is always taken and we have 50% of all predictions as mispredicted."

$ws.Range("E4").Value = "Main loop is always taken and we predicted it in all stages(except first time). Tested branch is always taken(0 == 0) and predictor is always true(except first time)"

$ws.Range("E5").Value = "Main loop is unconditional jump with NEGATIVE offset(backward) and it is always predicted. Tested branch is conditional jump with POSITIVE offset(forward) and is always mispredicted(50%)."
$ws.Range("E5").Characters(38, 8).Font.Bold = $true
$ws.Range("E5").Characters(131, 9).Font.Bold = $true

$ws.Range("E6").Value = "Main loop is always taken and predictor goes into a state of TAKEN after first time. Tested branch is conditional jump, but it is beqz zero(always true) and predictor becomes TAKEN after first time."

$ws.Range("E7").Value = "Main loop is always taken and predictor goes into a state of WEAKLY TAKEN after first time and then STRONGLY TAKEN. Tested branch is conditional jump, but it is beqz zero(always true) and predictor becomes WEAKLY TAKEN after first time and then STRONGLY TAKEN."

$ws.Range("E8").Value = "Main loop is always taken and history about this(pattern 11) is known for predictor after the first two times. Tested branch is bneq zero, label and it" + [char]0x2019 + "s always taken. History about this is known for predictor after the first two times(difference between this and other good predictors)"

# --- Row heights to fit the new comments ---
$ws.Rows.Item(3).RowHeight = 24
$ws.Rows.Item(4).RowHeight = 18
$ws.Rows.Item(5).RowHeight = 21
$ws.Rows.Item(6).RowHeight = 20.25
$ws.Rows.Item(7).RowHeight = 20.25
$ws.Rows.Item(8).RowHeight = 18

# --- Selection moves to E11 ---
[void]$ws.Range("E11").Select()
